$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Changes")

# Clear existing used range first so no stale cells remain
$ws.Cells.Clear()

$ws.Cells.Item(1, 1).Value = "path"
$ws.Cells.Item(1, 2).Value = "change"
$ws.Cells.Item(1, 3).Value = "file"

$ws.Cells.Item(2, 1).Value = "fieldPermissions.HIP_API_Transaction__c.API_Status__c"
$ws.Cells.Item(2, 2).Value = "Modified"
$ws.Cells.Item(2, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(3, 1).Value = "fieldPermissions.HIP_API_Transaction__c.API_Type__c"
$ws.Cells.Item(3, 2).Value = "Modified"
$ws.Cells.Item(3, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(4, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Date_and_Time_of_Submission_Failure__c"
$ws.Cells.Item(4, 2).Value = "Modified"
$ws.Cells.Item(4, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(5, 1).Value = "fieldPermissions.HIP_API_Transaction__c.DML_Operation__c"
$ws.Cells.Item(5, 2).Value = "Modified"
$ws.Cells.Item(5, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(6, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Destination__c"
$ws.Cells.Item(6, 2).Value = "Modified"
$ws.Cells.Item(6, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(7, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Endpoint__c"
$ws.Cells.Item(7, 2).Value = "Modified"
$ws.Cells.Item(7, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(8, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Error_Message__c"
$ws.Cells.Item(8, 2).Value = "Modified"
$ws.Cells.Item(8, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(9, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Exception_Details__c"
$ws.Cells.Item(9, 2).Value = "Modified"
$ws.Cells.Item(9, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(10, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Exception_Type__c"
$ws.Cells.Item(10, 2).Value = "Modified"
$ws.Cells.Item(10, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(11, 1).Value = "fieldPermissions.HIP_API_Transaction__c.External_Id_Field__c"
$ws.Cells.Item(11, 2).Value = "Modified"
$ws.Cells.Item(11, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(12, 1).Value = "fieldPermissions.HIP_API_Transaction__c.External_Id__c"
$ws.Cells.Item(12, 2).Value = "Modified"
$ws.Cells.Item(12, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(13, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Http_Method__c"
$ws.Cells.Item(13, 2).Value = "Modified"
$ws.Cells.Item(13, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(14, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Individual_Application__c"
$ws.Cells.Item(14, 2).Value = "Modified"
$ws.Cells.Item(14, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(15, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Is_Cached__c"
$ws.Cells.Item(15, 2).Value = "Modified"
$ws.Cells.Item(15, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(16, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Level__c"
$ws.Cells.Item(16, 2).Value = "Modified"
$ws.Cells.Item(16, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(17, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Log_Generated_Time__c"
$ws.Cells.Item(17, 2).Value = "Modified"
$ws.Cells.Item(17, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(18, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Message__c"
$ws.Cells.Item(18, 2).Value = "Modified"
$ws.Cells.Item(18, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(19, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Module__c"
$ws.Cells.Item(19, 2).Value = "Modified"
$ws.Cells.Item(19, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(20, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Object_Name__c"
$ws.Cells.Item(20, 2).Value = "Modified"
$ws.Cells.Item(20, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(21, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Origin__c"
$ws.Cells.Item(21, 2).Value = "Modified"
$ws.Cells.Item(21, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(22, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Parent_Request_Id__c"
$ws.Cells.Item(22, 2).Value = "Modified"
$ws.Cells.Item(22, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(23, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Record_Id__c"
$ws.Cells.Item(23, 2).Value = "Modified"
$ws.Cells.Item(23, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(24, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Request_Id__c"
$ws.Cells.Item(24, 2).Value = "Modified"
$ws.Cells.Item(24, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(25, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Request__c"
$ws.Cells.Item(25, 2).Value = "Modified"
$ws.Cells.Item(25, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(26, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Response__c"
$ws.Cells.Item(26, 2).Value = "Modified"
$ws.Cells.Item(26, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(27, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Running_User__c"
$ws.Cells.Item(27, 2).Value = "Modified"
$ws.Cells.Item(27, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(28, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Stack_Trace__c"
$ws.Cells.Item(28, 2).Value = "Modified"
$ws.Cells.Item(28, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(29, 1).Value = "fieldPermissions.HIP_API_Transaction__c.StatusCode__c"
$ws.Cells.Item(29, 2).Value = "Modified"
$ws.Cells.Item(29, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(30, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Status__c"
$ws.Cells.Item(30, 2).Value = "Modified"
$ws.Cells.Item(30, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(31, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Trace_Id__c"
$ws.Cells.Item(31, 2).Value = "Modified"
$ws.Cells.Item(31, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(32, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Transaction_End_Time__c"
$ws.Cells.Item(32, 2).Value = "Modified"
$ws.Cells.Item(32, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(33, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Transaction_Logs__c"
$ws.Cells.Item(33, 2).Value = "Modified"
$ws.Cells.Item(33, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(34, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Transaction_Start_Time__c"
$ws.Cells.Item(34, 2).Value = "Modified"
$ws.Cells.Item(34, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(35, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Transaction_Time_Spent_Milliseconds__c"
$ws.Cells.Item(35, 2).Value = "Modified"
$ws.Cells.Item(35, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(36, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Transaction_Type__c"
$ws.Cells.Item(36, 2).Value = "Modified"
$ws.Cells.Item(36, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(37, 1).Value = "fieldPermissions.HIP_API_Transaction__c.Type_Name__c"
$ws.Cells.Item(37, 2).Value = "Modified"
$ws.Cells.Item(37, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(38, 1).Value = "fieldPermissions.HIP_API_Transaction__c.UTC_End_Time_System__c"
$ws.Cells.Item(38, 2).Value = "Modified"
$ws.Cells.Item(38, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"

$ws.Cells.Item(39, 1).Value = "fieldPermissions.HIP_API_Transaction__c.UTC_Start_Time_System__c"
$ws.Cells.Item(39, 2).Value = "Modified"
$ws.Cells.Item(39, 3).Value = "Passport_API_Transactions_Read_Only.permissionset-meta.xml, TRV_API_Transactions_Read_Only.permissionset-meta.xml"


